$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" (D1) and "is_enabled" (E1) option columns entirely.
# Deleting the Range shifts the remaining cells (order_by, rem) left so
# they land on D1/E1, and the now-unused F1/G1 cells disappear.
$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
